$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.392.08'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.574.35'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3764'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3418'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.168'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07689'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.990'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.935'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001143'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = '1.573.73'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06727'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5280'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = '22.391.34'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.771'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '144.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.066'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").Value = '1.746.64'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.027'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.261'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08531'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02549'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06548'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.517'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.298'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6440'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.01%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.15%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6027'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.778'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.306'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("E51").Value = '  +2.52%  '
